# Insert a new data row at row 3 (pushing existing rows 3..40 down to 4..41)
# and populate it with a new price record, matching the commit's weekly update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 3 and below down by one row.
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with the same fixed values shared by every record in this
# sheet, plus the new date/price/origin data from the diff.
$ws.Cells.Item(3, 1).Value  = 10
$ws.Cells.Item(3, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(3, 3).Value  = "La Araucanía"
$ws.Cells.Item(3, 4).Value  = Get-Date -Year 2021 -Month 10 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(3, 5).Value  = 9
$ws.Cells.Item(3, 6).Value  = 100112026
$ws.Cells.Item(3, 7).Value  = "Haba"
$ws.Cells.Item(3, 8).Value  = "Sin especificar"
$ws.Cells.Item(3, 9).Value  = "Primera"
$ws.Cells.Item(3, 10).Value = 50
$ws.Cells.Item(3, 11).Value = 8000
$ws.Cells.Item(3, 12).Value = 9000
$ws.Cells.Item(3, 13).Value = 8600
$ws.Cells.Item(3, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(3, 15).Value = "Región Metropolitana"
$ws.Cells.Item(3, 16).Value = 344
$ws.Cells.Item(3, 17).Value = 25
$ws.Cells.Item(3, 18).Value = "Hortaliza"
